# Update cryptocurrency price/volume data per the Sun Apr 16 19:22:47 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" cell so numeric-looking text (e.g. "1.010", "0.5283")
# is stored as TEXT, matching the source data (every Price cell is an inline
# string in this sheet, never a real Number), instead of being auto-coerced by
# Excel into a Number/Date. The leading apostrophe forces text entry; ClearFormats
# then drops the resulting "quote prefix" cell style so formatting stays plain,
# matching the unstyled Price cells elsewhere in the sheet.
function Set-PriceText($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

Set-PriceText $ws.Cells.Item(2, 4) "30.753.52"
$ws.Cells.Item(2, 5).Value = "  +0.97%  "

Set-PriceText $ws.Cells.Item(3, 4) "2.151.58"
$ws.Cells.Item(3, 5).Value = "  +2.18%  "

Set-PriceText $ws.Cells.Item(4, 4) "1.010"
$ws.Cells.Item(4, 5).Value = "  +0.38%  "

$ws.Cells.Item(5, 5).Value = "  +5.55%  "

$ws.Cells.Item(6, 5).Value = "  +0.35%  "

Set-PriceText $ws.Cells.Item(7, 4) "0.5283"
$ws.Cells.Item(7, 5).Value = "  +1.11%  "

Set-PriceText $ws.Cells.Item(8, 4) "0.4572"
$ws.Cells.Item(8, 5).Value = "  +0.90%  "

Set-PriceText $ws.Cells.Item(9, 4) "54.34"
$ws.Cells.Item(9, 5).Value = "  +2.27%  "

Set-PriceText $ws.Cells.Item(10, 4) "0.09137"
$ws.Cells.Item(10, 5).Value = "  +2.52%  "

$ws.Cells.Item(11, 5).Value = "  +0.81%  "

Set-PriceText $ws.Cells.Item(12, 4) "24.97"
$ws.Cells.Item(12, 5).Value = "  +3.59%  "

Set-PriceText $ws.Cells.Item(13, 4) "2.143.68"
$ws.Cells.Item(13, 5).Value = "  +1.83%  "

Set-PriceText $ws.Cells.Item(14, 4) "6.928"
$ws.Cells.Item(14, 5).Value = "  +1.46%  "

Set-PriceText $ws.Cells.Item(15, 4) "8.175"
$ws.Cells.Item(15, 5).Value = "  +1.69%  "

Set-PriceText $ws.Cells.Item(16, 4) "102.52"
$ws.Cells.Item(16, 5).Value = "  +5.99%  "

Set-PriceText $ws.Cells.Item(17, 4) "0.00001179"
$ws.Cells.Item(17, 5).Value = "  +3.36%  "

$ws.Cells.Item(18, 5).Value = "  +0.29%  "

Set-PriceText $ws.Cells.Item(19, 4) "0.06728"
$ws.Cells.Item(19, 5).Value = "  +1.08%  "

Set-PriceText $ws.Cells.Item(20, 4) "19.63"
$ws.Cells.Item(20, 5).Value = "  +2.14%  "

$ws.Cells.Item(21, 5).Value = "  +0.37%  "

Set-PriceText $ws.Cells.Item(22, 4) "6.403"
$ws.Cells.Item(22, 5).Value = "  +1.02%  "

Set-PriceText $ws.Cells.Item(23, 4) "30.825.65"
$ws.Cells.Item(23, 5).Value = "  +1.00%  "

Set-PriceText $ws.Cells.Item(24, 4) "12.92"
$ws.Cells.Item(24, 5).Value = "  +3.67%  "

Set-PriceText $ws.Cells.Item(25, 4) "2.392"
$ws.Cells.Item(25, 5).Value = "  +1.55%  "

Set-PriceText $ws.Cells.Item(26, 4) "2.392.62"
$ws.Cells.Item(26, 5).Value = "  +1.74%  "

Set-PriceText $ws.Cells.Item(27, 4) "22.70"
$ws.Cells.Item(27, 5).Value = "  +2.11%  "

Set-PriceText $ws.Cells.Item(28, 4) "2.610"
$ws.Cells.Item(28, 5).Value = "  +3.25%  "

$ws.Cells.Item(29, 5).Value = "  +1.27%  "

Set-PriceText $ws.Cells.Item(30, 4) "137.03"
$ws.Cells.Item(30, 5).Value = "  +2.93%  "

Set-PriceText $ws.Cells.Item(31, 4) "1.223"
$ws.Cells.Item(31, 5).Value = "  +1.27%  "

Set-PriceText $ws.Cells.Item(32, 4) "0.1086"
$ws.Cells.Item(32, 5).Value = "  +1.29%  "

Set-PriceText $ws.Cells.Item(33, 4) "1.683"
$ws.Cells.Item(33, 5).Value = "  +1.65%  "

Set-PriceText $ws.Cells.Item(34, 4) "6.420"
$ws.Cells.Item(34, 5).Value = "  +0.19%  "

Set-PriceText $ws.Cells.Item(35, 4) "4.016"
$ws.Cells.Item(35, 5).Value = "  +1.81%  "

Set-PriceText $ws.Cells.Item(36, 4) "6.169"
$ws.Cells.Item(36, 5).Value = "  +6.26%  "

Set-PriceText $ws.Cells.Item(37, 4) "10.43"
$ws.Cells.Item(37, 5).Value = "  +0.21%  "

Set-PriceText $ws.Cells.Item(38, 4) "0.02663"
$ws.Cells.Item(38, 5).Value = "  +2.67%  "

Set-PriceText $ws.Cells.Item(39, 4) "0.06933"
$ws.Cells.Item(39, 5).Value = "  +1.30%  "

Set-PriceText $ws.Cells.Item(40, 4) "0.2342"
$ws.Cells.Item(40, 5).Value = "  +1.85%  "

Set-PriceText $ws.Cells.Item(41, 4) "12.69"
$ws.Cells.Item(41, 5).Value = "  -0.22%  "

Set-PriceText $ws.Cells.Item(42, 4) "0.6966"
$ws.Cells.Item(42, 5).Value = "  +1.45%  "

Set-PriceText $ws.Cells.Item(43, 4) "1.275"
$ws.Cells.Item(43, 5).Value = "  +2.06%  "

Set-PriceText $ws.Cells.Item(44, 4) "14.90"
$ws.Cells.Item(44, 5).Value = "  +6.34%  "

Set-PriceText $ws.Cells.Item(45, 4) "2.366"
$ws.Cells.Item(45, 5).Value = "  +2.35%  "

Set-PriceText $ws.Cells.Item(46, 4) "0.6489"
$ws.Cells.Item(46, 5).Value = "  +1.93%  "

Set-PriceText $ws.Cells.Item(49, 4) "1.261"
$ws.Cells.Item(49, 5).Value = "  +0.97%  "

# Rows 47/48 swap places in the ranking (BabyDogeCoin now ranks above PancakeSwap).
# Rows 50/51: WOONetwork drops out of the list, Aave moves up to 50, Cronos enters at 51.
$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-PriceText $ws.Cells.Item(47, 4) "0.00000000371"
$ws.Cells.Item(47, 5).Value = "  +4.84%  "

$ws.Cells.Item(48, 2).Value = "PancakeSwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-PriceText $ws.Cells.Item(48, 4) "3.760"
$ws.Cells.Item(48, 5).Value = "  +2.68%  "

$ws.Cells.Item(50, 2).Value = "Aave"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-PriceText $ws.Cells.Item(50, 4) "83.48"
$ws.Cells.Item(50, 5).Value = "  +0.03%  "

$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-PriceText $ws.Cells.Item(51, 4) "0.07336"
$ws.Cells.Item(51, 5).Value = "  +2.70%  "
